$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center the title across A1:E1 by merging the cells
# (the title's existing style already centers text horizontally)
$ws.Range("A1:E1").Merge() | Out-Null

# Format the date as literal text "07/01/2019" (not an actual date value),
# and set the client names, matching the exported report's layout
$ws.Range("A7").Value = "'07/01/2019"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = "John Doe"

$ws.Range("A8").Value = "'07/01/2019"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = "Jane Doe"

# Adjust column A width to fit new narrower content
# (nearest width Excel's pixel grid can represent to the target 11.609375)
$ws.Columns.Item(1).ColumnWidth = 10.8
